# Update the Point value for student 6410301026 (row 4, column C) from 38 to 41.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = 41
